$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Rows 161 <-> 162 are swapped in full (the running index in column A
#    stays put; every other column B..AC moves with its match).
# ---------------------------------------------------------------------------
function Swap-Rows($ws, $r1, $r2, $lastCol) {
    for ($c = 2; $c -le $lastCol; $c++) {
        $v1 = $ws.Cells.Item($r1, $c).Value
        $v2 = $ws.Cells.Item($r2, $c).Value
        $ws.Cells.Item($r1, $c).Value = $v2
        $ws.Cells.Item($r2, $c).Value = $v1
    }
}

Swap-Rows $ws 161 162 29
Swap-Rows $ws 166 167 29

# ---------------------------------------------------------------------------
# 2) Rows 168-171 keep their match but get refreshed closing odds.
# ---------------------------------------------------------------------------
$ws.Cells.Item(168, 14).Value = 3.1     # N168
$ws.Cells.Item(168, 16).Value = 2.15    # P168
$ws.Cells.Item(168, 18).Value = 1.95    # R168
$ws.Cells.Item(168, 19).Value = 1.95    # S168

$ws.Cells.Item(169, 18).Value = 1.91    # R169
$ws.Cells.Item(169, 19).Value = 1.99    # S169
$ws.Cells.Item(169, 21).Value = 2.025   # U169
$ws.Cells.Item(169, 22).Value = 1.825   # V169

$ws.Cells.Item(170, 14).Value = 3.25    # N170
$ws.Cells.Item(170, 15).Value = 3.25    # O170
$ws.Cells.Item(170, 16).Value = 2.25    # P170

$ws.Cells.Item(171, 18).Value = 1.83    # R171
$ws.Cells.Item(171, 19).Value = 2.07    # S171
$ws.Cells.Item(171, 21).Value = 1.925   # U171
$ws.Cells.Item(171, 22).Value = 1.925   # V171

# ---------------------------------------------------------------------------
# 3) Two brand-new fixtures are inserted as rows 172 and 173 (pushing the
#    former rows 172 and 173 down to 174 and 175). Row 171 has the exact
#    same "shape" (no FTHG/FTAG/FTR, no closing PL_Ah* columns) as these
#    four rows, so its formatting is cloned into 172:175 before the values
#    are written.
# ---------------------------------------------------------------------------
$ws.Range("A171:AC171").Copy()
$ws.Range("A172:AC175").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = $false

# --- new row 172 (id 170) ---
$ws.Cells.Item(172, 1).Value = 170
$ws.Cells.Item(172, 2).Value = 7609654
$ws.Cells.Item(172, 3).Value = "Sweden Allsvenskan"
$ws.Cells.Item(172, 4).Value = "Sweden Allsvenskan"
$ws.Cells.Item(172, 5).Value = 45389.47916666666
$ws.Cells.Item(172, 6).Value = "IFK Varnamo"
$ws.Cells.Item(172, 7).Value = "GAIS"
$ws.Cells.Item(172, 11).Value = 1.833
$ws.Cells.Item(172, 12).Value = 3.75
$ws.Cells.Item(172, 13).Value = 4
$ws.Cells.Item(172, 14).Value = 2.1
$ws.Cells.Item(172, 15).Value = 3.5
$ws.Cells.Item(172, 16).Value = 3.4
$ws.Cells.Item(172, 17).Value = -0.25
$ws.Cells.Item(172, 18).Value = 1.83
$ws.Cells.Item(172, 19).Value = 2.07
$ws.Cells.Item(172, 20).Value = 2.5
$ws.Cells.Item(172, 21).Value = 1.825
$ws.Cells.Item(172, 22).Value = 2.025
$ws.Cells.Item(172, 23).Value = 0
$ws.Cells.Item(172, 24).Value = 0
$ws.Cells.Item(172, 25).Value = 0
$ws.Cells.Item(172, 26).Value = 0
$ws.Cells.Item(172, 27).Value = 0

# --- new row 173 (id 171) ---
$ws.Cells.Item(173, 1).Value = 171
$ws.Cells.Item(173, 2).Value = 7609653
$ws.Cells.Item(173, 3).Value = "Sweden Allsvenskan"
$ws.Cells.Item(173, 4).Value = "Sweden Allsvenskan"
$ws.Cells.Item(173, 5).Value = 45389.47916666666
$ws.Cells.Item(173, 6).Value = "Vasteras SK"
$ws.Cells.Item(173, 7).Value = "Elfsborg"
$ws.Cells.Item(173, 11).Value = 3.8
$ws.Cells.Item(173, 12).Value = 3.6
$ws.Cells.Item(173, 13).Value = 1.909
$ws.Cells.Item(173, 14).Value = 3.4
$ws.Cells.Item(173, 15).Value = 3.6
$ws.Cells.Item(173, 16).Value = 2
$ws.Cells.Item(173, 17).Value = 0.5
$ws.Cells.Item(173, 18).Value = 1.86
$ws.Cells.Item(173, 19).Value = 2.04
$ws.Cells.Item(173, 20).Value = 2.75
$ws.Cells.Item(173, 21).Value = 1.85
$ws.Cells.Item(173, 22).Value = 2
$ws.Cells.Item(173, 23).Value = 0
$ws.Cells.Item(173, 24).Value = 0
$ws.Cells.Item(173, 25).Value = 0
$ws.Cells.Item(173, 26).Value = 0
$ws.Cells.Item(173, 27).Value = 0

# --- row 174 (was row 172, id 172, same match, recomputed closing odds) ---
$ws.Cells.Item(174, 1).Value = 172
$ws.Cells.Item(174, 2).Value = 7608281
$ws.Cells.Item(174, 3).Value = "Sweden Allsvenskan"
$ws.Cells.Item(174, 4).Value = "Sweden Allsvenskan"
$ws.Cells.Item(174, 5).Value = 45390.58333333334
$ws.Cells.Item(174, 6).Value = "Kalmar FF"
$ws.Cells.Item(174, 7).Value = "Sirius"
$ws.Cells.Item(174, 11).Value = 2.5
$ws.Cells.Item(174, 12).Value = 3.4
$ws.Cells.Item(174, 13).Value = 2.7
$ws.Cells.Item(174, 14).Value = 2.3
$ws.Cells.Item(174, 15).Value = 3.4
$ws.Cells.Item(174, 16).Value = 2.9
$ws.Cells.Item(174, 17).Value = -0.25
$ws.Cells.Item(174, 18).Value = 2.08
$ws.Cells.Item(174, 19).Value = 1.82
$ws.Cells.Item(174, 20).Value = 2.75
$ws.Cells.Item(174, 21).Value = 2
$ws.Cells.Item(174, 22).Value = 1.85
$ws.Cells.Item(174, 23).Value = 0
$ws.Cells.Item(174, 24).Value = 0
$ws.Cells.Item(174, 25).Value = 0
$ws.Cells.Item(174, 26).Value = 0
$ws.Cells.Item(174, 27).Value = 0

# --- row 175 (was row 173, id 173, same match, recomputed odds) ---
$ws.Cells.Item(175, 1).Value = 173
$ws.Cells.Item(175, 2).Value = 7607823
$ws.Cells.Item(175, 3).Value = "Sweden Allsvenskan"
$ws.Cells.Item(175, 4).Value = "Sweden Allsvenskan"
$ws.Cells.Item(175, 5).Value = 45390.59027777778
$ws.Cells.Item(175, 6).Value = "Djurgarden"
$ws.Cells.Item(175, 7).Value = "BK Hacken"
$ws.Cells.Item(175, 11).Value = 2.1
$ws.Cells.Item(175, 12).Value = 3.6
$ws.Cells.Item(175, 13).Value = 3.25
$ws.Cells.Item(175, 14).Value = 1.75
$ws.Cells.Item(175, 15).Value = 3.75
$ws.Cells.Item(175, 16).Value = 4.5
$ws.Cells.Item(175, 17).Value = -0.75
$ws.Cells.Item(175, 18).Value = 2.02
$ws.Cells.Item(175, 19).Value = 1.88
$ws.Cells.Item(175, 20).Value = 3
$ws.Cells.Item(175, 21).Value = 1.975
$ws.Cells.Item(175, 22).Value = 1.875
$ws.Cells.Item(175, 23).Value = 0
$ws.Cells.Item(175, 24).Value = 0
$ws.Cells.Item(175, 25).Value = 0
$ws.Cells.Item(175, 26).Value = 0
$ws.Cells.Item(175, 27).Value = 0
